$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename "PythonCode" -> "pythonCode"
# ---------------------------------------------------------------------------
$wsPython = $wb.Worksheets.Item(2)
$wsPython.Name = "pythonCode"

# ---------------------------------------------------------------------------
# 2. Restructure the "pythonCode" sheet: insert a new first column holding a
#    TestId. This shifts the previous column A (code) to B and column B
#    (result) to C.
# ---------------------------------------------------------------------------
$wsPython.Columns.Item(1).Insert()

# Grab the formatting of (old A4, now B4) before we touch it - this is the
# style used by the "Practice Qns" question rows we build below.
$wsPython.Range("B4").Copy()

# ---------------------------------------------------------------------------
# 3. Add the new "Practice Qns" sheet (right after "pythonCode") holding the
#    TestCaseId / pythonCode practice-question rows, and paste in the style
#    captured above for its question rows (B6:B13).
# ---------------------------------------------------------------------------
$wsPractice = $wb.Worksheets.Add($null, $wsPython)
$wsPractice.Name = "Practice Qns"
$wsPractice.Range("B6:B13").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4. Fill in "pythonCode" sheet content (TestId / code / Result columns).
# ---------------------------------------------------------------------------
$wsPython.Range("A1").Value = "TestId"
$wsPython.Range("A2").Value = "TC001"
$wsPython.Range("B2").Value = ""
$wsPython.Range("A3").Value = "TC002"
$wsPython.Range("B3").Value = "hello"
$wsPython.Range("A4").Value = "TC003"
$wsPython.Range("B4").Value = "print(""hello"")"
$wsPython.Range("B4").Style = "Normal"

# The code snippets that used to live in column B (rows 5-11) move away to
# the new "Practice Qns" sheet - remove them (and their formatting) here.
$wsPython.Range("B5:B11").Clear()

# Selection shown on this sheet
$wsPython.Range("B4").Select()

# ---------------------------------------------------------------------------
# 5. Fill in "Practice Qns" sheet content.
# ---------------------------------------------------------------------------
$wsPractice.Range("A1").Value = "TestCaseId"
$wsPractice.Range("B1").Value = "pythonCode"

$wsPractice.Range("A2").Value = "TC001"

$wsPractice.Range("A3").Value = "TC002"
$wsPractice.Range("B3").Value = "hello"

$wsPractice.Range("A4").Value = "TC003"
$wsPractice.Range("B4").Value = "print(""Hello"")"

$searchCode = "def search(input_list, num):`nif(num in input_list):`nprint(""Element Found"")`n\b`n\b`nelse:`nprint(""Not Found"")`n\b`n\b`n\b`n\b`nsearch([12, 23, 45, 67, 6, 90] , 12)"
$maxOnesCode = "def findMaxConsecutiveOnes(nums) :`ncount = 0`nresult = 0`nfor i in range(0, len(nums)):`nif (nums[i] == 0):`ncount = 0`n\b`n\b`nelse:`ncount+= 1`n\b`n\b`nresult = max(result, count)`n\b`n\b`nprint(result)`n\b`n\b`nfindMaxConsecutiveOnes([1,0,1,1,0,1])"
$findNumbersCode = "def findNumbers(nums):`nc=0`nfor i in nums:`nj=str(i)`nx=len(j)`nif x%2==0:`nc=c+1`n\b`n\b`n\b`n\b`nprint c`nreturn c`nfindNumbers([12,345,2,6,7896])"
$sortedSquaresCode = "def sortedSquares(nums):`nsquares_list = []`nfor i in range(0, len(nums)):`nsquare = nums[i] * nums[i];`nsquares_list.append(square)`n\b`n\b`nsorted_squares_list = sorted(squares_list)`nprint sorted_squares_list;`nreturn sorted_squares_list;`nsortedSquares([-7,-3,2,3,11])"

$wsPractice.Range("B6").Value = $searchCode
$wsPractice.Range("B7").Value = $searchCode
$wsPractice.Range("B8").Value = $maxOnesCode
$wsPractice.Range("B9").Value = $maxOnesCode
$wsPractice.Range("B10").Value = $findNumbersCode
$wsPractice.Range("B11").Value = $findNumbersCode
$wsPractice.Range("B12").Value = $sortedSquaresCode
$wsPractice.Range("B13").Value = $sortedSquaresCode

# Drop the auto-grown row heights that typing multi-line text creates.
$wsPractice.Range("6:13").EntireRow.AutoFit()

$wsPractice.Columns.Item(1).ColumnWidth = 10.57
$wsPractice.Columns.Item(2).ColumnWidth = 243.43

$wsPractice.Range("A2").Select()

# ---------------------------------------------------------------------------
# 6. Re-activate "pythonCode" as the active sheet (matches activeTab=1).
# ---------------------------------------------------------------------------
$wsPython.Activate()

Write-Host "Edit complete"
